$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "x_Starch_Day_sp_exchange"
$ws.Range("B2").Value = -4.186128843675311
$ws.Range("C2").Value = -999999.8532632083
$ws.Range("D2").Value = -3.711215751192143
$ws.Range("E2").Value = $true

# Delete row 3 entirely (shift rows up)
$ws.Rows.Item(3).Delete()
